$d = $word.ActiveDocument

# --- Change 1: paragraph "Definir quem sera responsavel por cada parte do teste."
#     becomes "Criacao dos cenarios de teste em BDD -> Jhody Mike" and is followed
#     by two new paragraphs describing further test-execution scenarios. ---
$p1 = $d.Paragraphs.Item(55)
$r1 = $p1.Range

$xmlFrag1 = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:lang w:val='pt-BR'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t>Criação dos cenários de teste em BDD -&gt; Jhody Mike</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val='pt-BR'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t>Execução dos cenários de teste (funcionalidades login de usuários, suporte de usuários) -&gt; Jhody Mike</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val='pt-BR'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:lastRenderedPageBreak/><w:t>Execução dos cenários de teste (funcionalidades visualização de estoque e adição de novo produto, visualização de dashboard) Fernando</w:t></w:r></w:p>
"@
$r1.InsertXML($xmlFrag1)

# --- Change 2: paragraph "De acordo com a analise do QA responsavel pelos testes..."
#     loses its lastRenderedPageBreak (moved above) and its text is split into
#     several runs: "do" -> "dos" and "QA" -> "QAs", with proofErr spell-check
#     markers wrapping the inserted "QAs". ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*De acordo com a an*lise*") {
        $p2 = $p
        break
    }
}
$r2 = $p2.Range

$xmlFrag2 = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:lang w:val='pt-BR'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t>De acordo com a análise do</w:t></w:r><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t>QA</w:t></w:r><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t>s</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:rPr><w:lang w:val='pt-BR'/></w:rPr><w:t xml:space='preserve'> responsável pelos testes estima-se que cada execução de cenários de teste levará cerca de 5 a 10 minutos, para executar, evidenciar e reportar caso tenha falhas.</w:t></w:r></w:p>
"@
$r2.InsertXML($xmlFrag2)
